$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.ClearFormats()
}

# Row 2
Set-TextValue "D2" '29.296.81'
Set-TextValue "E2" '  +0.68%  '

# Row 3
Set-TextValue "D3" '1.844.03'

# Row 4
Set-TextValue "D4" '0.9993'
Set-TextValue "E4" '  -0.07%  '

# Row 5
Set-TextValue "D5" '241.10'
Set-TextValue "E5" '  -0.11%  '

# Row 6
Set-TextValue "D6" '0.6726'
Set-TextValue "E6" '  -1.81%  '

# Row 8
Set-TextValue "D8" '0.07444'
Set-TextValue "E8" '  +0.13%  '

# Row 9
Set-TextValue "D9" '0.2944'
Set-TextValue "E9" '  -2.45%  '

# Row 10
Set-TextValue "D10" '22.94'
Set-TextValue "E10" '  -0.63%  '

# Row 11
Set-TextValue "E11" '  +0.63%  '

# Row 12
Set-TextValue "D12" '1.837.82'
Set-TextValue "E12" '  -0.25%  '

# Row 13
Set-TextValue "E13" '  -0.89%  '

# Row 14
Set-TextValue "D14" '0.6719'
Set-TextValue "E14" '  -1.60%  '

# Row 15
Set-TextValue "D15" '85.83'
Set-TextValue "E15" '  -1.81%  '

# Row 16
Set-TextValue "D16" '6.148'
Set-TextValue "E16" '  -0.32%  '

# Row 17
Set-TextValue "D17" '29.284.09'
Set-TextValue "E17" '  +0.62%  '

# Row 18
Set-TextValue "D18" '0.000008318'
Set-TextValue "E18" '  +2.17%  '

# Row 19
Set-TextValue "D19" '229.06'

# Row 20
Set-TextValue "D20" '12.53'

# Row 22
Set-TextValue "D22" '7.176'
Set-TextValue "E22" '  -2.79%  '

# Row 23
Set-TextValue "E23" '  -0.03%  '

# Row 24
Set-TextValue "D24" '161.02'
Set-TextValue "E24" '  +0.52%  '

# Row 25
Set-TextValue "D25" '8.710'
Set-TextValue "E25" '  -0.66%  '

# Row 26
Set-TextValue "D26" '0.1406'
Set-TextValue "E26" '  -3.03%  '

# Row 27
Set-TextValue "D27" '18.04'
Set-TextValue "E27" '  -0.27%  '

# Row 28
Set-TextValue "D28" '1.517'
Set-TextValue "E28" '  +0.22%  '

# Row 29
Set-TextValue "D29" '4.161'
Set-TextValue "E29" '  -2.45%  '

# Row 30
Set-TextValue "E30" '  -1.47%  '

# Row 31
Set-TextValue "D31" '1.193'
Set-TextValue "E31" '  -0.07%  '

# Row 32
Set-TextValue "D32" '0.05306'
Set-TextValue "E32" '  +1.40%  '

# Row 33
Set-TextValue "B33" 'LidoDAOToken'
Set-TextValue "C33" 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
Set-TextValue "D33" '1.877'
Set-TextValue "E33" '  +1.57%  '

# Row 34
Set-TextValue "B34" 'ImmutableX'
Set-TextValue "C34" 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextValue "D34" '0.7534'
Set-TextValue "E34" '  -0.57%  '

# Row 35
Set-TextValue "D35" '1.138'
Set-TextValue "E35" '  +0.42%  '

# Row 36
Set-TextValue "D36" '2.678'
Set-TextValue "E36" '  -0.41%  '

# Row 37
Set-TextValue "D37" '1.320.86'
Set-TextValue "E37" '  +1.12%  '

# Row 38
Set-TextValue "D38" '0.01807'
Set-TextValue "E38" '  -1.61%  '

# Row 39
Set-TextValue "D39" '2.726'
Set-TextValue "E39" '  +0.00%  '

# Row 40
Set-TextValue "D40" '0.9205'
Set-TextValue "E40" '  -1.25%  '

# Row 41
Set-TextValue "D41" '5.976'
Set-TextValue "E41" '  +2.78%  '

# Row 42
Set-TextValue "D42" '0.08351'
Set-TextValue "E42" '  +13.60%  '

# Row 43
Set-TextValue "E43" '  +0.28%  '

# Row 44
Set-TextValue "D44" '102.85'
Set-TextValue "E44" '  -1.87%  '

# Row 45
Set-TextValue "B45" 'BabyDogeCoin'
Set-TextValue "C45" 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
Set-TextValue "D45" '0.00000000125'
Set-TextValue "E45" '  +1.71%  '

# Row 46
Set-TextValue "B46" 'RocketPoolETH'
Set-TextValue "C46" 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
Set-TextValue "D46" '1.986.84'
Set-TextValue "E46" '  +0.23%  '

# Row 47
Set-TextValue "E47" '  -0.56%  '

# Row 48
Set-TextValue "D48" '1.776'
Set-TextValue "E48" '  +0.37%  '

# Row 49
Set-TextValue "D49" '64.05'
Set-TextValue "E49" '  -1.06%  '

# Row 50
Set-TextValue "D50" '9.139'
Set-TextValue "E50" '  -3.92%  '

# Row 51
Set-TextValue "D51" '0.05951'
Set-TextValue "E51" '  +0.18%  '
